$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.324.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.709.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.78%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5296'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.30%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2641'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.38%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06550'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.97'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.86%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07635'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.90%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.569'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.34%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.732.99'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.47%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.945.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.70%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5745'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.04%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8195'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.44%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.21%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.308.67'
$ws.Range('D18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.10'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.75%  '

$ws.Range('E20').Value = '  +0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.681'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.69%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.966'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.53%  '

$ws.Range('E26').Value = '  +7.84%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1219'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.77%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.268'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.33%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.49%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05391'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.294'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.53%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.493'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.86%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.415'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.95%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.640'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.868'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.60%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9518'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.56%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.417'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.92%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5872'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01629'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.868'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.07%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.041.20'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.09%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8389'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.03%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.853.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.57%  '

$ws.Range('E46').Value = '  +4.99%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.59%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4496'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.98%  '

$ws.Range('E49').Value = '  +0.47%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.063'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.99%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05241'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.78%  '
